$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-8 from 46081 to 46082 (serial date 2026-02-28 -> 2026-03-01)
$ws.Range("C2:C8").Value = 46082
